$d = $word.ActiveDocument

# Merge the split runs "{{ tecnico_" + "2" + " }}" into a single run "{{ tecnico_2 }}"
$d.Content.Find.Execute("{{ tecnico_2 }}", $false, $false, $false, $false, $false, $true, 1, $false, "{{ tecnico_2 }}", 2)

# Merge the split runs "{{ h_" + "2" + " }}" into a single run "{{ h_2 }}"
$d.Content.Find.Execute("{{ h_2 }}", $false, $false, $false, $false, $false, $true, 1, $false, "{{ h_2 }}", 2)
